# Update the cryptocurrency price/volume figures on the active sheet.
# Only column D (Price) and column E (Volume(1h)) change, per the diff;
# columns A-C (index, coin name, link) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = new price (or $null if unchanged); E = new volume string }
$updates = @{
    2  = @{ D = "67.253.37";  E = "  -1.98%  " }
    3  = @{ D = "2.665.62";   E = "  -1.37%  " }
    4  = @{ D = $null;        E = "  +0.01%  " }
    5  = @{ D = "597.94";     E = "  -0.20%  " }
    6  = @{ D = "165.62";     E = "  +3.16%  " }
    7  = @{ D = $null;        E = "  +0.02%  " }
    8  = @{ D = $null;        E = "  +0.17%  " }
    9  = @{ D = "2.665.01";   E = "  -1.36%  " }
    10 = @{ D = $null;        E = "  +0.89%  " }
    11 = @{ D = $null;        E = "  +1.20%  " }
    12 = @{ D = $null;        E = "  -1.01%  " }
    13 = @{ D = $null;        E = "  -1.62%  " }
    14 = @{ D = "27.72";      E = "  -2.20%  " }
    15 = @{ D = "3.151.56";   E = "  -1.28%  " }
    16 = @{ D = $null;        E = "  -2.70%  " }
    17 = @{ D = "67.165.17";  E = "  -1.98%  " }
    18 = @{ D = "2.670.72";   E = "  -0.97%  " }
    19 = @{ D = "11.70";      E = "  -1.32%  " }
    20 = @{ D = "7.63";       E = "  -0.12%  " }
    21 = @{ D = "362.93";     E = "  -0.85%  " }
    22 = @{ D = $null;        E = "  -3.71%  " }
    23 = @{ D = $null;        E = "  -2.25%  " }
    24 = @{ D = $null;        E = "  -4.65%  " }
    26 = @{ D = "70.61";      E = "  -5.21%  " }
    27 = @{ D = "10.01";      E = "  +0.69%  " }
    28 = @{ D = "2.800.04";   E = "  -1.36%  " }
    30 = @{ D = $null;        E = "  +0.00%  " }
    31 = @{ D = "554.72";     E = "  -4.32%  " }
    32 = @{ D = "7.98";       E = "  -3.24%  " }
    33 = @{ D = $null;        E = "  -4.43%  " }
    34 = @{ D = $null;        E = "  -1.43%  " }
    35 = @{ D = "0.130";      E = "  -2.33%  " }
    36 = @{ D = $null;        E = "  +0.00%  " }
    37 = @{ D = $null;        E = "  -5.83%  " }
    38 = @{ D = "19.50";      E = "  -1.74%  " }
    39 = @{ D = "156.04";     E = "  -3.53%  " }
    40 = @{ D = $null;        E = "  -2.10%  " }
    41 = @{ D = "5.29";       E = "  -2.19%  " }
    42 = @{ D = $null;        E = "  -4.64%  " }
    43 = @{ D = $null;        E = "  +0.44%  " }
    44 = @{ D = $null;        E = "  +0.03%  " }
    45 = @{ D = "2.51";       E = "  -7.09%  " }
    46 = @{ D = $null;        E = "  -0.79%  " }
    47 = @{ D = $null;        E = "  -6.75%  " }
    48 = @{ D = $null;        E = "  -2.13%  " }
    49 = @{ D = "152.67";     E = "  -3.28%  " }
    50 = @{ D = "3.83";       E = "  -3.11%  " }
    51 = @{ D = $null;        E = "  -3.41%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        # Column D holds plain text (e.g. "19.50", "2.664.25"); force text
        # formatting before assigning so Excel doesn't reinterpret numeric-
        # looking values as numbers and drop significant trailing zeros.
        $dCell = $ws.Range("D$row")
        $dCell.NumberFormat = "@"
        $dCell.Value = $vals.D
    }
    $ws.Range("E$row").Value = $vals.E
}
